$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set row height + values + formatting + hyperlinks for the newly purchased BOM lines
$dataFont = 'Helvetica Neue'
$dataFontSize = 14
$dataFontColor = 3355443
$currencyFormat = '"$"#,##0.00_);[Red]("$"#,##0.00)'

# Row 53: C-PD001-600
$r = 53
$ws.Rows.Item($r).RowHeight = 18
$ws.Cells.Item($r, 1).Value = 1
$ws.Cells.Item($r, 2).Value = 1
$ws.Cells.Item($r, 3).Value = 'EA'
$ws.Cells.Item($r, 4).Value = 'C-PD001-600'
$ws.Cells.Item($r, 5).Value = 'Capacitor - 600V, 716P Series, Polypropylene, Capacitance: .001 uF'
$ws.Cells.Item($r, 6).Value = 1.1
$ws.Cells.Item($r, 7).Value = 1.1
$ws.Range("A$r`:C$r").Font.Name = $dataFont
$ws.Range("A$r`:C$r").Font.Size = $dataFontSize
$ws.Range("A$r`:C$r").Font.Color = $dataFontColor
$ws.Range("E$r").Font.Name = $dataFont
$ws.Range("E$r").Font.Size = $dataFontSize
$ws.Range("E$r").Font.Color = $dataFontColor
$ws.Range("F$r`:G$r").Font.Name = $dataFont
$ws.Range("F$r`:G$r").Font.Size = $dataFontSize
$ws.Range("F$r`:G$r").Font.Color = $dataFontColor
$ws.Range("F$r`:G$r").NumberFormat = $currencyFormat
$ws.Hyperlinks.Add($ws.Range("D$r"), 'https://www.amplifiedparts.com/products/capacitor-600v-716p-series-polypropylene', '', '', 'C-PD001-600') | Out-Null
$ws.Range("D$r").Style = 'Hyperlink'

# Row 54: C-LD022-630
$r = 54
$ws.Rows.Item($r).RowHeight = 18
$ws.Cells.Item($r, 1).Value = 1
$ws.Cells.Item($r, 2).Value = 1
$ws.Cells.Item($r, 3).Value = 'EA'
$ws.Cells.Item($r, 4).Value = 'C-LD022-630'
$ws.Cells.Item($r, 5).Value = 'Capacitor - 630V, Polypropylene, radial leads, Capacitance: .022 uF'
$ws.Cells.Item($r, 6).Value = 0.36
$ws.Cells.Item($r, 7).Value = 0.36
$ws.Range("A$r`:C$r").Font.Name = $dataFont
$ws.Range("A$r`:C$r").Font.Size = $dataFontSize
$ws.Range("A$r`:C$r").Font.Color = $dataFontColor
$ws.Range("E$r").Font.Name = $dataFont
$ws.Range("E$r").Font.Size = $dataFontSize
$ws.Range("E$r").Font.Color = $dataFontColor
$ws.Range("F$r`:G$r").Font.Name = $dataFont
$ws.Range("F$r`:G$r").Font.Size = $dataFontSize
$ws.Range("F$r`:G$r").Font.Color = $dataFontColor
$ws.Range("F$r`:G$r").NumberFormat = $currencyFormat
$ws.Hyperlinks.Add($ws.Range("D$r"), 'https://www.amplifiedparts.com/products/capacitor-630v-polypropylene-radial-leads', '', '', 'C-LD022-630') | Out-Null
$ws.Range("D$r").Style = 'Hyperlink'

# Row 55: R-B33K
$r = 55
$ws.Rows.Item($r).RowHeight = 18
$ws.Cells.Item($r, 1).Value = 5
$ws.Cells.Item($r, 2).Value = 5
$ws.Cells.Item($r, 3).Value = 'PKG/5'
$ws.Cells.Item($r, 4).Value = 'R-B33K'
$ws.Cells.Item($r, 5).Value = 'Resistors - 1W, Carbon Film, Resistance: 33 kOhm'
$ws.Cells.Item($r, 6).Value = 2
$ws.Cells.Item($r, 7).Value = 2
$ws.Range("A$r`:C$r").Font.Name = $dataFont
$ws.Range("A$r`:C$r").Font.Size = $dataFontSize
$ws.Range("A$r`:C$r").Font.Color = $dataFontColor
$ws.Range("E$r").Font.Name = $dataFont
$ws.Range("E$r").Font.Size = $dataFontSize
$ws.Range("E$r").Font.Color = $dataFontColor
$ws.Range("F$r`:G$r").Font.Name = $dataFont
$ws.Range("F$r`:G$r").Font.Size = $dataFontSize
$ws.Range("F$r`:G$r").Font.Color = $dataFontColor
$ws.Range("F$r`:G$r").NumberFormat = $currencyFormat
$ws.Hyperlinks.Add($ws.Range("D$r"), 'https://www.amplifiedparts.com/products/resistors-1w-carbon-film', '', '', 'R-B33K') | Out-Null
$ws.Range("D$r").Style = 'Hyperlink'

# Row 56: P-H1590B
$r = 56
$ws.Rows.Item($r).RowHeight = 18
$ws.Cells.Item($r, 1).Value = 1
$ws.Cells.Item($r, 2).Value = 1
$ws.Cells.Item($r, 3).Value = 'EA'
$ws.Cells.Item($r, 4).Value = 'P-H1590B'
$ws.Cells.Item($r, 5).Value = 'Chassis Box - Hammond, Unpainted Aluminum, 4.37" x 2.37" x 1.22"'
$ws.Cells.Item($r, 6).Value = 7.07
$ws.Cells.Item($r, 7).Value = 7.07
$ws.Range("A$r`:C$r").Font.Name = $dataFont
$ws.Range("A$r`:C$r").Font.Size = $dataFontSize
$ws.Range("A$r`:C$r").Font.Color = $dataFontColor
$ws.Range("E$r").Font.Name = $dataFont
$ws.Range("E$r").Font.Size = $dataFontSize
$ws.Range("E$r").Font.Color = $dataFontColor
$ws.Range("F$r`:G$r").Font.Name = $dataFont
$ws.Range("F$r`:G$r").Font.Size = $dataFontSize
$ws.Range("F$r`:G$r").Font.Color = $dataFontColor
$ws.Range("F$r`:G$r").NumberFormat = $currencyFormat
$ws.Hyperlinks.Add($ws.Range("D$r"), 'https://www.amplifiedparts.com/products/chassis-box-hammond-unpainted-aluminum-437-x-237-x-122', '', '', 'P-H1590B') | Out-Null
$ws.Range("D$r").Style = 'Hyperlink'

# Row 57: P-K302
$r = 57
$ws.Rows.Item($r).RowHeight = 18
$ws.Cells.Item($r, 1).Value = 2
$ws.Cells.Item($r, 2).Value = 2
$ws.Cells.Item($r, 3).Value = 'EA'
$ws.Cells.Item($r, 4).Value = 'P-K302'
$ws.Cells.Item($r, 5).Value = 'Knob - Chicken Head, Raised, Color: Black'
$ws.Cells.Item($r, 6).Value = 1.4
$ws.Cells.Item($r, 7).Value = 2.8
$ws.Range("A$r`:C$r").Font.Name = $dataFont
$ws.Range("A$r`:C$r").Font.Size = $dataFontSize
$ws.Range("A$r`:C$r").Font.Color = $dataFontColor
$ws.Range("E$r").Font.Name = $dataFont
$ws.Range("E$r").Font.Size = $dataFontSize
$ws.Range("E$r").Font.Color = $dataFontColor
$ws.Range("F$r`:G$r").Font.Name = $dataFont
$ws.Range("F$r`:G$r").Font.Size = $dataFontSize
$ws.Range("F$r`:G$r").Font.Color = $dataFontColor
$ws.Range("F$r`:G$r").NumberFormat = $currencyFormat
$ws.Hyperlinks.Add($ws.Range("D$r"), 'https://www.amplifiedparts.com/products/knob-chicken-head-raised', '', '', 'P-K302') | Out-Null
$ws.Range("D$r").Style = 'Hyperlink'

# Row 58: C-LD001-630
$r = 58
$ws.Rows.Item($r).RowHeight = 18
$ws.Cells.Item($r, 1).Value = 2
$ws.Cells.Item($r, 2).Value = 2
$ws.Cells.Item($r, 3).Value = 'EA'
$ws.Cells.Item($r, 4).Value = 'C-LD001-630'
$ws.Cells.Item($r, 5).Value = 'Capacitor - 630V, Polypropylene, radial leads, Capacitance: .001 uF'
$ws.Cells.Item($r, 6).Value = 0.34
$ws.Cells.Item($r, 7).Value = 0.68
$ws.Range("A$r`:C$r").Font.Name = $dataFont
$ws.Range("A$r`:C$r").Font.Size = $dataFontSize
$ws.Range("A$r`:C$r").Font.Color = $dataFontColor
$ws.Range("E$r").Font.Name = $dataFont
$ws.Range("E$r").Font.Size = $dataFontSize
$ws.Range("E$r").Font.Color = $dataFontColor
$ws.Range("F$r`:G$r").Font.Name = $dataFont
$ws.Range("F$r`:G$r").Font.Size = $dataFontSize
$ws.Range("F$r`:G$r").Font.Color = $dataFontColor
$ws.Range("F$r`:G$r").NumberFormat = $currencyFormat
$ws.Hyperlinks.Add($ws.Range("D$r"), 'https://www.amplifiedparts.com/products/capacitor-630v-polypropylene-radial-leads', '', '', 'C-LD001-630') | Out-Null
$ws.Range("D$r").Style = 'Hyperlink'

# Row 59: S-HLW6
$r = 59
$ws.Rows.Item($r).RowHeight = 18
$ws.Cells.Item($r, 1).Value = 25
$ws.Cells.Item($r, 2).Value = 25
$ws.Cells.Item($r, 3).Value = 'PKG/5'
$ws.Cells.Item($r, 4).Value = 'S-HLW6'
$ws.Cells.Item($r, 5).Value = 'Washer - Internal Tooth Lock, Zinc, Size: #6'
$ws.Cells.Item($r, 6).Value = 0.45
$ws.Cells.Item($r, 7).Value = 2.25
$ws.Range("A$r`:C$r").Font.Name = $dataFont
$ws.Range("A$r`:C$r").Font.Size = $dataFontSize
$ws.Range("A$r`:C$r").Font.Color = $dataFontColor
$ws.Range("E$r").Font.Name = $dataFont
$ws.Range("E$r").Font.Size = $dataFontSize
$ws.Range("E$r").Font.Color = $dataFontColor
$ws.Range("F$r`:G$r").Font.Name = $dataFont
$ws.Range("F$r`:G$r").Font.Size = $dataFontSize
$ws.Range("F$r`:G$r").Font.Color = $dataFontColor
$ws.Range("F$r`:G$r").NumberFormat = $currencyFormat
$ws.Hyperlinks.Add($ws.Range("D$r"), 'https://www.amplifiedparts.com/products/washer-internal-tooth-lock-zinc', '', '', 'S-HLW6') | Out-Null
$ws.Range("D$r").Style = 'Hyperlink'

# Row 60: P-H395
$r = 60
$ws.Rows.Item($r).RowHeight = 18
$ws.Cells.Item($r, 1).Value = 2
$ws.Cells.Item($r, 2).Value = 2
$ws.Cells.Item($r, 3).Value = 'EA'
$ws.Cells.Item($r, 4).Value = 'P-H395'
$ws.Cells.Item($r, 5).Value = 'Switch - Rotary, 3 Poles, 3 Position'
$ws.Cells.Item($r, 6).Value = 2.5
$ws.Cells.Item($r, 7).Value = 5
$ws.Range("A$r`:C$r").Font.Name = $dataFont
$ws.Range("A$r`:C$r").Font.Size = $dataFontSize
$ws.Range("A$r`:C$r").Font.Color = $dataFontColor
$ws.Range("E$r").Font.Name = $dataFont
$ws.Range("E$r").Font.Size = $dataFontSize
$ws.Range("E$r").Font.Color = $dataFontColor
$ws.Range("F$r`:G$r").Font.Name = $dataFont
$ws.Range("F$r`:G$r").Font.Size = $dataFontSize
$ws.Range("F$r`:G$r").Font.Color = $dataFontColor
$ws.Range("F$r`:G$r").NumberFormat = $currencyFormat
$ws.Hyperlinks.Add($ws.Range("D$r"), 'https://www.amplifiedparts.com/products/switch-rotary-3-poles-3-position', '', '', 'P-H395') | Out-Null
$ws.Range("D$r").Style = 'Hyperlink'

# Row 61: S-HLW38
$r = 61
$ws.Rows.Item($r).RowHeight = 18
$ws.Cells.Item($r, 1).Value = 15
$ws.Cells.Item($r, 2).Value = 15
$ws.Cells.Item($r, 3).Value = 'PKG/5'
$ws.Cells.Item($r, 4).Value = 'S-HLW38'
$ws.Cells.Item($r, 5).Value = 'Washer - Internal Tooth Lock, Zinc, Size: 3/8"'
$ws.Cells.Item($r, 6).Value = 0.65
$ws.Cells.Item($r, 7).Value = 1.95
$ws.Range("A$r`:C$r").Font.Name = $dataFont
$ws.Range("A$r`:C$r").Font.Size = $dataFontSize
$ws.Range("A$r`:C$r").Font.Color = $dataFontColor
$ws.Range("E$r").Font.Name = $dataFont
$ws.Range("E$r").Font.Size = $dataFontSize
$ws.Range("E$r").Font.Color = $dataFontColor
$ws.Range("F$r`:G$r").Font.Name = $dataFont
$ws.Range("F$r`:G$r").Font.Size = $dataFontSize
$ws.Range("F$r`:G$r").Font.Color = $dataFontColor
$ws.Range("F$r`:G$r").NumberFormat = $currencyFormat
$ws.Hyperlinks.Add($ws.Range("D$r"), 'https://www.amplifiedparts.com/products/washer-internal-tooth-lock-zinc', '', '', 'S-HLW38') | Out-Null
$ws.Range("D$r").Style = 'Hyperlink'

# Row 62: S-HS632-38
$r = 62
$ws.Rows.Item($r).RowHeight = 18
$ws.Cells.Item($r, 1).Value = 25
$ws.Cells.Item($r, 2).Value = 25
$ws.Cells.Item($r, 3).Value = 'PKG/5'
$ws.Cells.Item($r, 4).Value = 'S-HS632-38'
$ws.Cells.Item($r, 5).Value = 'Screw - 6/32, Phillips, Pan Head, Machine, Zinc, Length: 3/8"'
$ws.Cells.Item($r, 6).Value = 0.4
$ws.Cells.Item($r, 7).Value = 2
$ws.Range("A$r`:C$r").Font.Name = $dataFont
$ws.Range("A$r`:C$r").Font.Size = $dataFontSize
$ws.Range("A$r`:C$r").Font.Color = $dataFontColor
$ws.Range("E$r").Font.Name = $dataFont
$ws.Range("E$r").Font.Size = $dataFontSize
$ws.Range("E$r").Font.Color = $dataFontColor
$ws.Range("F$r`:G$r").Font.Name = $dataFont
$ws.Range("F$r`:G$r").Font.Size = $dataFontSize
$ws.Range("F$r`:G$r").Font.Color = $dataFontColor
$ws.Range("F$r`:G$r").NumberFormat = $currencyFormat
$ws.Hyperlinks.Add($ws.Range("D$r"), 'https://www.amplifiedparts.com/products/screw-632-phillips-pan-head-machine-zinc', '', '', 'S-HS632-38') | Out-Null
$ws.Range("D$r").Style = 'Hyperlink'

# Row 63: S-H173
$r = 63
$ws.Rows.Item($r).RowHeight = 18
$ws.Cells.Item($r, 1).Value = 20
$ws.Cells.Item($r, 2).Value = 20
$ws.Cells.Item($r, 3).Value = 'PKG/4'
$ws.Cells.Item($r, 4).Value = 'S-H173'
$ws.Cells.Item($r, 5).Value = 'Standoff - #6-32, Female, Aluminum, Length: 3/4"'
$ws.Cells.Item($r, 6).Value = 2.2
$ws.Cells.Item($r, 7).Value = 11
$ws.Range("A$r`:C$r").Font.Name = $dataFont
$ws.Range("A$r`:C$r").Font.Size = $dataFontSize
$ws.Range("A$r`:C$r").Font.Color = $dataFontColor
$ws.Range("E$r").Font.Name = $dataFont
$ws.Range("E$r").Font.Size = $dataFontSize
$ws.Range("E$r").Font.Color = $dataFontColor
$ws.Range("F$r`:G$r").Font.Name = $dataFont
$ws.Range("F$r`:G$r").Font.Size = $dataFontSize
$ws.Range("F$r`:G$r").Font.Color = $dataFontColor
$ws.Range("F$r`:G$r").NumberFormat = $currencyFormat
$ws.Hyperlinks.Add($ws.Range("D$r"), 'https://www.amplifiedparts.com/products/standoff-6-32-female-aluminum', '', '', 'S-H173') | Out-Null
$ws.Range("D$r").Style = 'Hyperlink'

# Row 64: S-H172
$r = 64
$ws.Rows.Item($r).RowHeight = 18
$ws.Cells.Item($r, 1).Value = 20
$ws.Cells.Item($r, 2).Value = 20
$ws.Cells.Item($r, 3).Value = 'PKG/4'
$ws.Cells.Item($r, 4).Value = 'S-H172'
$ws.Cells.Item($r, 5).Value = 'Standoff - #6-32, Female, Aluminum, Length: 1/2"'
$ws.Cells.Item($r, 6).Value = 2
$ws.Cells.Item($r, 7).Value = 10
$ws.Range("A$r`:C$r").Font.Name = $dataFont
$ws.Range("A$r`:C$r").Font.Size = $dataFontSize
$ws.Range("A$r`:C$r").Font.Color = $dataFontColor
$ws.Range("E$r").Font.Name = $dataFont
$ws.Range("E$r").Font.Size = $dataFontSize
$ws.Range("E$r").Font.Color = $dataFontColor
$ws.Range("F$r`:G$r").Font.Name = $dataFont
$ws.Range("F$r`:G$r").Font.Size = $dataFontSize
$ws.Range("F$r`:G$r").Font.Color = $dataFontColor
$ws.Range("F$r`:G$r").NumberFormat = $currencyFormat
$ws.Hyperlinks.Add($ws.Range("D$r"), 'https://www.amplifiedparts.com/products/standoff-6-32-female-aluminum', '', '', 'S-H172') | Out-Null
$ws.Range("D$r").Style = 'Hyperlink'

# Row 65: T-5751-PS-TAD
$r = 65
$ws.Rows.Item($r).RowHeight = 18
$ws.Cells.Item($r, 1).Value = 1
$ws.Cells.Item($r, 2).Value = 1
$ws.Cells.Item($r, 3).Value = 'EA'
$ws.Cells.Item($r, 4).Value = 'T-5751-PS-TAD'
$ws.Cells.Item($r, 5).Value = 'Vacuum Tube - 5751, Tube Amp Doctor, Premium Selected'
$ws.Cells.Item($r, 6).Value = 18.95
$ws.Cells.Item($r, 7).Value = 18.95
$ws.Range("A$r`:C$r").Font.Name = $dataFont
$ws.Range("A$r`:C$r").Font.Size = $dataFontSize
$ws.Range("A$r`:C$r").Font.Color = $dataFontColor
$ws.Range("E$r").Font.Name = $dataFont
$ws.Range("E$r").Font.Size = $dataFontSize
$ws.Range("E$r").Font.Color = $dataFontColor
$ws.Range("F$r`:G$r").Font.Name = $dataFont
$ws.Range("F$r`:G$r").Font.Size = $dataFontSize
$ws.Range("F$r`:G$r").Font.Color = $dataFontColor
$ws.Range("F$r`:G$r").NumberFormat = $currencyFormat
$ws.Hyperlinks.Add($ws.Range("D$r"), 'https://www.amplifiedparts.com/products/vacuum-tube-5751-tube-amp-doctor-premium-selected', '', '', 'T-5751-PS-TAD') | Out-Null
$ws.Range("D$r").Style = 'Hyperlink'

# Row 66: T-7025-HG-TAD
$r = 66
$ws.Rows.Item($r).RowHeight = 18
$ws.Cells.Item($r, 1).Value = 6
$ws.Cells.Item($r, 2).Value = 6
$ws.Cells.Item($r, 3).Value = 'EA'
$ws.Cells.Item($r, 4).Value = 'T-7025-HG-TAD'
$ws.Cells.Item($r, 5).Value = 'Vacuum Tube - 7025, Tube Amp Doctor, High Grade'
$ws.Cells.Item($r, 6).Value = 19.95
$ws.Cells.Item($r, 7).Value = 119.7
$ws.Range("A$r`:C$r").Font.Name = $dataFont
$ws.Range("A$r`:C$r").Font.Size = $dataFontSize
$ws.Range("A$r`:C$r").Font.Color = $dataFontColor
$ws.Range("E$r").Font.Name = $dataFont
$ws.Range("E$r").Font.Size = $dataFontSize
$ws.Range("E$r").Font.Color = $dataFontColor
$ws.Range("F$r`:G$r").Font.Name = $dataFont
$ws.Range("F$r`:G$r").Font.Size = $dataFontSize
$ws.Range("F$r`:G$r").Font.Color = $dataFontColor
$ws.Range("F$r`:G$r").NumberFormat = $currencyFormat
$ws.Hyperlinks.Add($ws.Range("D$r"), 'https://www.amplifiedparts.com/products/vacuum-tube-7025-tube-amp-doctor-high-grade', '', '', 'T-7025-HG-TAD') | Out-Null
$ws.Range("D$r").Style = 'Hyperlink'

# Row 67: T-6L6WGC-TAD-MQ
$r = 67
$ws.Rows.Item($r).RowHeight = 18
$ws.Cells.Item($r, 1).Value = 1
$ws.Cells.Item($r, 2).Value = 1
$ws.Cells.Item($r, 3).Value = 'EA'
$ws.Cells.Item($r, 4).Value = 'T-6L6WGC-TAD-MQ'
$ws.Cells.Item($r, 5).Value = 'Vacuum Tube - 6L6WGC, Tube Amp Doctor, Single or Matched: Matched Quad'
$ws.Cells.Item($r, 6).Value = 98
$ws.Cells.Item($r, 7).Value = 98
$ws.Range("A$r`:C$r").Font.Name = $dataFont
$ws.Range("A$r`:C$r").Font.Size = $dataFontSize
$ws.Range("A$r`:C$r").Font.Color = $dataFontColor
$ws.Range("E$r").Font.Name = $dataFont
$ws.Range("E$r").Font.Size = $dataFontSize
$ws.Range("E$r").Font.Color = $dataFontColor
$ws.Range("F$r`:G$r").Font.Name = $dataFont
$ws.Range("F$r`:G$r").Font.Size = $dataFontSize
$ws.Range("F$r`:G$r").Font.Color = $dataFontColor
$ws.Range("F$r`:G$r").NumberFormat = $currencyFormat
$ws.Hyperlinks.Add($ws.Range("D$r"), 'https://www.amplifiedparts.com/products/vacuum-tube-6l6wgc-tube-amp-doctor', '', '', 'T-6L6WGC-TAD-MQ') | Out-Null
$ws.Range("D$r").Style = 'Hyperlink'

# Update the view: scroll down and select the newly added rows (matches the author's last save)
$win = $excel.ActiveWindow
$win.ScrollRow = 35
$win.ScrollColumn = 1
$ws.Range("A53:G67").Select()

Write-Output "Added 15 new BOM rows (53-67) from amplifiedparts.com"
